# Evaporation Runs from last week
# - Rename the sheet "Run0" -> "Trace2"
# - Move the active selection from G24 -> H29 (last worked cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Run0")
$ws.Name = "Trace2"

$ws.Activate()
$ws.Range("H29").Select()
